$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -6.091599999999995
$ws.Range("B9").Value = 8.528400000000003
$ws.Range("B18").Value = 4.749400000000005
$ws.Range("B20").Value = 5.585599999999999
